$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.025.84'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.40%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.303.09'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.89%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.58'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.41'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.63%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.57%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.608'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.94'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.04%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0911'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.47'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.51%  '

$ws.Range("E13").Value = '  +1.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.978'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.42%  '

$ws.Range("E15").Value = '  -3.45%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.649.99'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.95%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.287.89'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.094.47'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.50%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.72'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.68%  '

$ws.Range("E20").Value = '  -0.32%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.20'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.56'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '259.37'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.65%  '

$ws.Range("E24").Value = '  -0.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.98'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +4.35%  '

$ws.Range("E26").Value = '  +0.60%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.95'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.88'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.92%  '

$ws.Range("E29").Value = '  +0.65%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.97'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.60%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '164.41'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -5.96%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0886'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.91'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.52%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.88'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.121'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.38%  '

$ws.Range("E36").Value = '  +0.60%  '

$ws.Range("E37").Value = '  +0.82%  '

$ws.Range("E38").Value = '  +9.77%  '

$ws.Range("E39").Value = '  -2.12%  '

$ws.Range("E40").Value = '  -3.54%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '101.74'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +14.28%  '

$ws.Range("E42").Value = '  +1.23%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '71.02'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.34%  '

$ws.Range("E44").Value = '  -1.96%  '

$ws.Range("E45").Value = '  -0.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.15'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '114.38'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.81%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '79.02'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +8.02%  '

$ws.Range("E49").Value = '  -0.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.33'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.14%  '

$ws.Range("E51").Value = '  +2.16%  '
